$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 'aa'
$ws.Range("J2").Value = 'Agree/Accept'
$ws.Range("I5").Value = 'ba'
$ws.Range("J5").Value = 'Appreciation'
$ws.Range("I16").Value = 'sd'
$ws.Range("J16").Value = 'Statement-non-opinion'
$ws.Range("I19").Value = 'sd'
$ws.Range("J19").Value = 'Statement-non-opinion'
$ws.Range("I21").Value = '%'
$ws.Range("J21").Value = 'Uninterpretable'
$ws.Range("I26").Value = '%'
$ws.Range("J26").Value = 'Uninterpretable'
$ws.Range("I36").Value = 'sd'
$ws.Range("J36").Value = 'Statement-non-opinion'
$ws.Range("I44").Value = 'sd'
$ws.Range("J44").Value = 'Statement-non-opinion'
$ws.Range("I54").Value = 'sd'
$ws.Range("J54").Value = 'Statement-non-opinion'
$ws.Range("I55").Value = 'sd'
$ws.Range("J55").Value = 'Statement-non-opinion'
$ws.Range("I69").Value = 'aa'
$ws.Range("J69").Value = 'Agree/Accept'
$ws.Range("I94").Value = 'ba'
$ws.Range("J94").Value = 'Appreciation'
$ws.Range("I95").Value = 'sv'
$ws.Range("J95").Value = 'Statement-opinion'
$ws.Range("I96").Value = 'sv'
$ws.Range("J96").Value = 'Statement-opinion'
$ws.Range("I97").Value = 'sv'
$ws.Range("J97").Value = 'Statement-opinion'
$ws.Range("I102").Value = 'ba'
$ws.Range("J102").Value = 'Appreciation'
$ws.Range("I105").Value = 'qy'
$ws.Range("J105").Value = 'Yes-No-Question'
$ws.Range("I106").Value = 'sv'
$ws.Range("J106").Value = 'Statement-opinion'
$ws.Range("I112").Value = 'ba'
$ws.Range("J112").Value = 'Appreciation'
$ws.Range("I113").Value = 'sd'
$ws.Range("J113").Value = 'Statement-non-opinion'
$ws.Range("I117").Value = 'sd'
$ws.Range("J117").Value = 'Statement-non-opinion'
$ws.Range("I118").Value = 'ba'
$ws.Range("J118").Value = 'Appreciation'
$ws.Range("I119").Value = 'sd'
$ws.Range("J119").Value = 'Statement-non-opinion'
$ws.Range("I130").Value = 'ba'
$ws.Range("J130").Value = 'Appreciation'
$ws.Range("I153").Value = 'aa'
$ws.Range("J153").Value = 'Agree/Accept'
$ws.Range("I163").Value = 'aa'
$ws.Range("J163").Value = 'Agree/Accept'
$ws.Range("I192").Value = 'ba'
$ws.Range("J192").Value = 'Appreciation'
$ws.Range("I193").Value = 'aa'
$ws.Range("J193").Value = 'Agree/Accept'
$ws.Range("I202").Value = 'sd'
$ws.Range("J202").Value = 'Statement-non-opinion'
$ws.Range("I205").Value = 'ba'
$ws.Range("J205").Value = 'Appreciation'
$ws.Range("I211").Value = 'sd'
$ws.Range("J211").Value = 'Statement-non-opinion'
$ws.Range("I236").Value = 'sd'
$ws.Range("J236").Value = 'Statement-non-opinion'
$ws.Range("I239").Value = 'sv'
$ws.Range("J239").Value = 'Statement-opinion'
$ws.Range("I241").Value = 'sd'
$ws.Range("J241").Value = 'Statement-non-opinion'
$ws.Range("I246").Value = 'sd'
$ws.Range("J246").Value = 'Statement-non-opinion'
